# Apply the edits described by the diff:
# 1) Insert a new row at row 5 (shifts old row5 "Natural Gas" -> row6,
#    and old row6 "EIA Carbon..." (H only) -> row7)
# 2) Populate new row 5 with "Coal" data
# 3) Populate new row 7 with "Natural Gas - DAC" data (H7 already carries
#    the shifted "EIA Carbon Dioxide Emission Coefficiencts" text)
# 4) Add row 12 with "Direct Air Capture" in A12
# 5) Widen column A to fit the new longer labels
# 6) Update the active cell selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Insert a new row above row 5 ---
$ws.Rows.Item(5).Insert()

# --- 2) Fill the new row 5 with Coal data, mirroring the Natural Gas row pattern ---
$ws.Range("A5").Value = "Coal"
$ws.Range("B5").Value = -90
$ws.Range("C5").Value = 93.3
$ws.Range("B5:D5").Interior.Color = 65535
$ws.Range("E5").Formula = "=C5*0.000001/0.00000105505585262"
$ws.Range("F5").Formula = "=E5*B5/100"

# --- 3) Fill the new row 7 with Natural Gas - DAC data ---
$ws.Range("A7").Value = "Natural Gas - DAC"
$ws.Range("B7").Value = -98
$ws.Range("C7").Value = 53.07
$ws.Range("B7:D7").Interior.Color = 65535
$ws.Range("E7").Formula = "=C7*0.000001/0.00000105505585262"
$ws.Range("F7").Formula = "=E7*B7/100"

# --- 4) Add Direct Air Capture label in A12 ---
$ws.Range("A12").Value = "Direct Air Capture"

# --- 5) Widen column A for the new, longer row labels ---
$ws.Columns.Item(1).ColumnWidth = 15.3

# --- 6) Update selection to match the saved view state ---
$null = $ws.Range("A14").Select()
